# Update "want to go" counts (column F) on three worksheets to reflect
# newly generated output (gh-pages data refresh at commit 456a3b4).
#
# Sheet "展览"   (Exhibitions)  - sheet index 1
# Sheet "演出"   (Performances) - sheet index 2
# Sheet "全部类型" (All types)   - sheet index 4
# Sheet "本地生活" (Local life)  - unchanged, not touched here

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (exhibitions) ---
$wsExhibition.Range("F4").Value = 316
$wsExhibition.Range("F5").Value = 1273
$wsExhibition.Range("F7").Value = 306
$wsExhibition.Range("F8").Value = 1108
$wsExhibition.Range("F10").Value = 6931
$wsExhibition.Range("F14").Value = 7824
$wsExhibition.Range("F16").Value = 48
$wsExhibition.Range("F17").Value = 5434
$wsExhibition.Range("F19").Value = 2305
$wsExhibition.Range("F21").Value = 4541
$wsExhibition.Range("F22").Value = 266
$wsExhibition.Range("F26").Value = 305
$wsExhibition.Range("F28").Value = 6
$wsExhibition.Range("F29").Value = 2054
$wsExhibition.Range("F31").Value = 232
$wsExhibition.Range("F33").Value = 25
$wsExhibition.Range("F34").Value = 540
$wsExhibition.Range("F37").Value = 1397
$wsExhibition.Range("F38").Value = 23
$wsExhibition.Range("F40").Value = 2130

# --- Sheet "演出" (performances) ---
$wsPerformance.Range("F4").Value = 31

# --- Sheet "全部类型" (all types, aggregate of all sheets) ---
$wsAll.Range("F7").Value = 316
$wsAll.Range("F8").Value = 1273
$wsAll.Range("F10").Value = 306
$wsAll.Range("F11").Value = 1108
$wsAll.Range("F13").Value = 6931
$wsAll.Range("F17").Value = 7824
$wsAll.Range("F19").Value = 48
$wsAll.Range("F20").Value = 5435
$wsAll.Range("F22").Value = 2305
$wsAll.Range("F24").Value = 4541
$wsAll.Range("F25").Value = 266
$wsAll.Range("F30").Value = 31
$wsAll.Range("F31").Value = 305
$wsAll.Range("F33").Value = 2054
$wsAll.Range("F35").Value = 232
$wsAll.Range("F37").Value = 25
$wsAll.Range("F38").Value = 540
$wsAll.Range("F42").Value = 1397
$wsAll.Range("F43").Value = 23
$wsAll.Range("F45").Value = 2130
